$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("02-11-2021", 4.23, 3.34, 3.42),
    @("03-11-2021", 4.21, 3.38, 3.56),
    @("04-11-2021", 4.21, 3.43, 3.7),
    @("05-11-2021", 4.11, 3.4, 3.68)
)

$startRow = 213
$endRow = $startRow + $data.Length - 1

# Pre-format column A for these new rows as text so the dd-mm-yyyy style
# date strings (e.g. "02-11-2021") are kept as literal text, matching the
# existing "Serie" column entries (which are stored as plain shared
# strings, not date serials).
$rangeAddr = "A" + $startRow + ":A" + $endRow
$ws.Range($rangeAddr).NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $entry[0]
    # Drop back to the default (unstyled) look now that the text is safely
    # stored as a string, so the new rows match the plain formatting of
    # the rest of the "Serie" column.
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
